$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Clear the two numeric "id" values in column A (rows 46-47), keep styles ---
$ws.Range("A46").ClearContents()
$ws.Range("A47").ClearContents()

# --- Fill in the new English/Portuguese phrase rows (D/E, rows 89-102) ---
# Order below matches the exact order the cells were authored in, so new
# shared-string entries land in the same sequence as the source workbook.

# E91 then E92 get a thin box border drawn around them (E91 keeps its
# pre-existing highlight fill, E92 has no fill) - applying the borders here,
# in this order, reproduces the two new cellXfs entries in the right slots.
$ws.Range("E91").Borders.LineStyle = 1
$ws.Range("E92").Borders.LineStyle = 1

# Row 89 reuses two already-existing shared strings.
$ws.Range("D89").Value = "How is this possible (ráuu is diz póssibol) "
$ws.Range("E89").Value = "Como isto é possivel"

# Row 90
$ws.Range("E90").Value = "Como vão as coisas"
$ws.Range("D90").Value = "How's it going? (Rauzit Going?)"

# Row 91 (D only for now; E91 is filled in further below)
$ws.Range("D91").Value = "How are you doing? "

# Row 92 (D only for now; E92 is filled in further below)
$ws.Range("D92").Value = "How are You"

# Row 93
$ws.Range("D93").Value = "What do you mean? ( Uóriu mim? )"
$ws.Range("E93").Value = "O que você quer dizer?"

# Row 94
$ws.Range("E94").Value = "Que horas são?"
$ws.Range("D94").Value = "What time is it? (Uótaimizit)"

# Row 92 - E92 gets the boxed "como vai?" text (border already applied above)
$ws.Range("E92").Value = "como vai?"

# Row 95
$ws.Range("D95").Value = "Are you busy?"
$ws.Range("E95").Value = "Você está ocupado?"

# Row 96
$ws.Range("D96").Value = "long time no see!"
$ws.Range("E96").Value = "Quanto tempo não te vejo!"

# Row 97
$ws.Range("D97").Value = "Where are you from? "
$ws.Range("E97").Value = "De onde você é?"

# Row 98
$ws.Range("E98").Value = "Quantos anos você tem?"
$ws.Range("D98").Value = "How old are you? (Rall OLDAR YOU?)"

# Row 99
$ws.Range("E99").Value = "Prazer em conhecê-lo!"
$ws.Range("D99").Value = "Nice to meet you! (nice to meetchiu)"

# Row 100
$ws.Range("E100").Value = "Tenha um bom dia"
$ws.Range("D100").Value = "Have a nice day ( REVA nice day)"

# Row 101
$ws.Range("D101").Value = "What´s your name? (Uótisíor name"

# Row 102 (brand-new row)
$ws.Range("D102").Value = "How can I help you!"
$ws.Range("E102").Value = "Como posso ajudar você"

# Row 91 - E91 gets the same boxed phrase as E92 (reuses the shared string);
# its border was already applied above, on top of the pre-existing highlight fill.
$ws.Range("E91").Value = "como vai?"

# --- Update selection to match the end state of the edit ---
$ws.Range("D102").Select()
